$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.001, 20.17678),
    @(0.001995199519951995, 65.08519605027043),
    @(0.003878587858785879, 76.89611597274758),
    @(0.006061606160616062, 81.27920348918057),
    @(0.01001030103010301, 85.06747124243778),
    @(0.0196947694769477, 91.26910925345638),
    @(0.02503460346034603, 93.08187038278433),
    @(0.02728182818281828, 95.15994164090588),
    @(0.02959325932593259, 95.38835776820723),
    @(0.03447294729472947, 97.38877792660372),
    @(0.03650615061506151, 97.5585755635397),
    @(0.04028362836283628, 98.94595242443549),
    @(0.04210281028102811, 98.60474807513995),
    @(0.04409320932093209, 99.71119918119514),
    @(0.0444998499849985, 99.73557863411435),
    @(0.04550575057505751, 99.65853333491032),
    @(0.04599799979998, 99.69050082413528),
    @(0.04993599359935994, 100.7848440152658),
    @(0.05041754175417542, 100.8236718949981),
    @(0.05148764876487649, 100.7671325757923),
    @(0.05196919691969197, 100.790716261659),
    @(0.05486918691869187, 101.6063703001341),
    @(0.05513671367136714, 101.630110916602),
    @(0.05572527252725273, 101.6063280976701),
    @(0.05598209820982099, 101.6208375196905),
    @(0.05704150415041504, 101.9041965987892),
    @(0.05774777477747775, 101.8832742196092),
    @(0.05797249724972497, 101.895493002133),
    @(0.05895699569956996, 102.1719098610917),
    @(0.05987728772877288, 102.1372505394933),
    @(0.06158945894589459, 103.2747684551664),
    @(0.06323742374237423, 102.9766878617648),
    @(0.06429682968296831, 103.0109135197394),
    @(0.06501380138013801, 103.003851995792),
    @(0.06583778377837785, 103.0076057867868),
    @(0.0665975597559756, 102.9810696454496),
    @(0.06692929292929294, 102.9954069944486),
    @(0.07093149314931493, 104.0976997498527),
    @(0.0713916391639164, 104.1320595746894),
    @(0.0722905290529053, 104.1009616066534),
    @(0.07362816281628164, 104.1179268829378),
    @(0.07454845484548456, 104.0771647836916),
    @(0.0748908890889089, 104.0916018540722),
    @(0.0801879187918792, 105.2312465002684),
    @(0.08136503650365037, 105.2061658209206),
    @(0.08234953495349535, 105.2143295030023),
    @(0.08333403340334034, 105.2121540331565),
    @(0.08771077107710772, 105.2131172473107),
    @(0.0886952695269527, 105.2106821947993),
    @(0.08966906690669067, 105.2199445890384),
    @(0.0906963696369637, 105.1857805448556),
    @(0.09248344834483449, 105.5426353086867),
    @(0.0942919291929193, 105.1729214668576),
    @(0.09575797579757976, 105.2913194568232),
    @(0.09862586258625863, 103.7267255675627),
    @(0.1007125712571257, 100.5782170946842),
    @(0.108, 76.77531)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

